$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 67.77251700000001
$ws.Range("H2").Value = 203.317551
$ws.Range("I2").Value = 0.4079637943863715
$ws.Range("J2").Value = 0.4079637943863715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.215523666666666
$ws.Range("N2").Value = 18.646571
$ws.Range("O2").Value = 0.2852115546146347
$ws.Range("P2").Value = 0.2852115546146347
$ws.Range("Q2").Value = 421.241683363069
$ws.Range("R2").Value = 3791.175150267621
$ws.Range("S2").Value = 0.1163559880234222
$ws.Range("T2").Value = 0.1163559880234222

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 67.77251700000001
$ws.Range("H3").Value = 203.317551
$ws.Range("I3").Value = 0.4079637943863715
$ws.Range("J3").Value = 0.4079637943863715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.484070666666666
$ws.Range("N3").Value = 28.452212
$ws.Range("O3").Value = 0.4351952762116512
$ws.Range("P3").Value = 0.4351952762116512
$ws.Range("Q3").Value = 642.7593404858681
$ws.Range("R3").Value = 5784.834064372812
$ws.Range("S3").Value = 0.1775439161823302
$ws.Range("T3").Value = 0.1775439161823302

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 67.77251700000001
$ws.Range("H4").Value = 203.317551
$ws.Range("I4").Value = 0.4079637943863715
$ws.Range("J4").Value = 0.4079637943863715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.093084
$ws.Range("N4").Value = 18.279252
$ws.Range("O4").Value = 0.2795931691737141
$ws.Range("P4").Value = 0.2795931691737141
$ws.Range("Q4").Value = 412.9436389724281
$ws.Range("R4").Value = 3716.492750751852
$ws.Range("S4").Value = 0.1140638901806191
$ws.Range("T4").Value = 0.1140638901806191

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 60.97760633333333
$ws.Range("H5").Value = 182.932819
$ws.Range("I5").Value = 0.3670611149405164
$ws.Range("J5").Value = 0.3670611149405164
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.215523666666666
$ws.Range("N5").Value = 18.646571
$ws.Range("O5").Value = 0.2852115546146347
$ws.Range("P5").Value = 0.2852115546146347
$ws.Range("Q5").Value = 379.0077553015165
$ws.Range("R5").Value = 3411.069797713649
$ws.Range("S5").Value = 0.1046900712307658
$ws.Range("T5").Value = 0.1046900712307658

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 60.97760633333333
$ws.Range("H6").Value = 182.932819
$ws.Range("I6").Value = 0.3670611149405164
$ws.Range("J6").Value = 0.3670611149405164
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.484070666666666
$ws.Range("N6").Value = 28.452212
$ws.Range("O6").Value = 0.4351952762116512
$ws.Range("P6").Value = 0.4351952762116512
$ws.Range("Q6").Value = 578.3159275495142
$ws.Range("R6").Value = 5204.843347945628
$ws.Range("S6").Value = 0.1597432633030947
$ws.Range("T6").Value = 0.1597432633030947

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 60.97760633333333
$ws.Range("H7").Value = 182.932819
$ws.Range("I7").Value = 0.3670611149405164
$ws.Range("J7").Value = 0.3670611149405164
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.093084
$ws.Range("N7").Value = 18.279252
$ws.Range("O7").Value = 0.2795931691737141
$ws.Range("P7").Value = 0.2795931691737141
$ws.Range("Q7").Value = 371.541677507932
$ws.Range("R7").Value = 3343.875097571388
$ws.Range("S7").Value = 0.1026277804066559
$ws.Range("T7").Value = 0.1026277804066559

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.37372866666667
$ws.Range("H8").Value = 112.121186
$ws.Range("I8").Value = 0.2249750906731122
$ws.Range("J8").Value = 0.2249750906731122
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.215523666666666
$ws.Range("N8").Value = 18.646571
$ws.Range("O8").Value = 0.2852115546146347
$ws.Range("P8").Value = 0.2852115546146347
$ws.Range("Q8").Value = 232.2972950392451
$ws.Range("R8").Value = 2090.675655353206
$ws.Range("S8").Value = 0.06416549536044673
$ws.Range("T8").Value = 0.06416549536044673

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.37372866666667
$ws.Range("H9").Value = 112.121186
$ws.Range("I9").Value = 0.2249750906731122
$ws.Range("J9").Value = 0.2249750906731122
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.484070666666666
$ws.Range("N9").Value = 28.452212
$ws.Range("O9").Value = 0.4351952762116512
$ws.Range("P9").Value = 0.4351952762116512
$ws.Range("Q9").Value = 354.4550837514925
$ws.Range("R9").Value = 3190.095753763432
$ws.Range("S9").Value = 0.09790809672622634
$ws.Range("T9").Value = 0.09790809672622634

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.37372866666667
$ws.Range("H10").Value = 112.121186
$ws.Range("I10").Value = 0.2249750906731122
$ws.Range("J10").Value = 0.2249750906731122
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.093084
$ws.Range("N10").Value = 18.279252
$ws.Range("O10").Value = 0.2795931691737141
$ws.Range("P10").Value = 0.2795931691737141
$ws.Range("Q10").Value = 227.7212681592081
$ws.Range("R10").Value = 2049.491413432872
$ws.Range("S10").Value = 0.06290149858643913
$ws.Range("T10").Value = 0.06290149858643912
